# Update the TPM-derived values in the LR-pairs worksheet (Insl5-Rxfp4)
# as produced by re-running the analysis scripts with the new TPM data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 5.002662333333333
$ws.Range("N2").Value = 15.007987
$ws.Range("O2").Value = 0.3527593062265809
$ws.Range("P2").Value = 0.3527593062265809
$ws.Range("Q2").Value = 0.6526023040456668
$ws.Range("R2").Value = 5.873420736411
$ws.Range("S2").Value = 0.3527593062265809
$ws.Range("T2").Value = 0.3527593062265809

# Row 3
$ws.Range("O3").Value = 0.3625979570169652
$ws.Range("P3").Value = 0.3625979570169652
$ws.Range("S3").Value = 0.3625979570169652
$ws.Range("T3").Value = 0.3625979570169652

# Row 4
$ws.Range("O4").Value = 0.2846427367564539
$ws.Range("P4").Value = 0.2846427367564539
$ws.Range("S4").Value = 0.2846427367564539
$ws.Range("T4").Value = 0.2846427367564539
